$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-18 from 45188 to 45189
$ws.Range("C2:C18").Value = 45189
